$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: wipe all existing cell content/formatting from the old table ---
$ws.Cells.Clear()
$ws.Rows.Item(1).ClearFormats()
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(3).ClearFormats()

# --- Header row (new column order: Vessel Type, Barcode, UMI Length, Spacer Length, Location) ---
$ws.Range("A1").Value = "Vessel Type"
$ws.Range("B1").Value = "Barcode"
$ws.Range("C1").Value = "UMI Length"
$ws.Range("D1").Value = "Spacer Length"
$ws.Range("E1").Value = "Location"

# --- Data rows ---
$ws.Range("A2").Value = "Eppendorf96"
$ws.Range("B2").Value = 12345
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "Inline First Read"

$ws.Range("A3").Value = "Eppendorf96"
$ws.Range("B3").Value = 34567
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Before Second Index Read"

$ws.Range("A4").Value = "Eppendorf96"
$ws.Range("B4").Value = 66789
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Inline Second Read"

$ws.Range("A5").Value = "Eppendorf96"
$ws.Range("B5").Value = 77891
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "Inline Second Read"

$ws.Range("A6").Value = "MatrixTube075"
$ws.Range("B6").Value = 87654
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Before First Read"

$ws.Range("A7").Value = "MatrixTube075"
$ws.Range("B7").Value = 87654
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Before Second Read"

$ws.Range("A8").Value = "Eppendorf96"
$ws.Range("B8").Value = 77891
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Before First Read"

# --- Formatting ---
# Whole header row is bold Arial 10, dark grey (FF222222)
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").Font.Name = "Arial"
$ws.Range("A1:E1").Font.Size = 10
$ws.Range("A1:E1").Font.Color = 2236962

# A1 alone is the same Arial/size/color family but NOT bold
$ws.Range("A1").Font.Bold = $false

# C3:D3 use the bold (default Calibri 11) font
$ws.Range("C3:D3").Font.Bold = $true

# --- Selection / view state ---
$ws.Range("E8").Select()
